# Apply weekly fruit/vegetable price update by rewriting the shuffled rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values per row for columns D (date serial), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen), P (Precio $/Kg). Row 6 is unchanged and omitted.
$rows = @(
    @{ Row = 2;  D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 3;  D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 },
    @{ Row = 4;  D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 },
    @{ Row = 5;  D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 7;  D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 8;  D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 },
    @{ Row = 9;  D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 },
    @{ Row = 10; D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 },
    @{ Row = 11; D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 },
    @{ Row = 12; D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 },
    @{ Row = 13; D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
}
